# Add a new translation text entry (row 18) to the "Translation" sheet,
# mirroring the existing SingleUseId17 row but with a new unique id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B18").Value = "SingleUseId18"
$ws.Range("C18").Value = "Default"
$ws.Range("D18").Value = "Center"
$ws.Range("E18").Value = "LTR"
$ws.Range("F18").Value = "New Text"
